$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I18").Value = "ba"
$ws.Range("J18").Value = "Appreciation"
$ws.Range("I25").Value = "sv"
$ws.Range("J25").Value = "Statement-opinion"
$ws.Range("I46").Value = "b"
$ws.Range("J46").Value = "Acknowledge (Backchannel)"
$ws.Range("I51").Value = "ba"
$ws.Range("J51").Value = "Appreciation"
$ws.Range("I57").Value = "sv"
$ws.Range("J57").Value = "Statement-opinion"
$ws.Range("I63").Value = "ba"
$ws.Range("J63").Value = "Appreciation"
$ws.Range("I70").Value = "b"
$ws.Range("J70").Value = "Acknowledge (Backchannel)"
$ws.Range("I74").Value = "b"
$ws.Range("J74").Value = "Acknowledge (Backchannel)"
$ws.Range("I83").Value = "b"
$ws.Range("J83").Value = "Acknowledge (Backchannel)"
$ws.Range("I109").Value = "ba"
$ws.Range("J109").Value = "Appreciation"
$ws.Range("I110").Value = "b"
$ws.Range("J110").Value = "Acknowledge (Backchannel)"
$ws.Range("I114").Value = "ba"
$ws.Range("J114").Value = "Appreciation"
$ws.Range("I117").Value = "ba"
$ws.Range("J117").Value = "Appreciation"
$ws.Range("I137").Value = "aa"
$ws.Range("J137").Value = "Agree/Accept"
$ws.Range("I138").Value = "aa"
$ws.Range("J138").Value = "Agree/Accept"
$ws.Range("I143").Value = "aa"
$ws.Range("J143").Value = "Agree/Accept"
$ws.Range("I147").Value = "sv"
$ws.Range("J147").Value = "Statement-opinion"
$ws.Range("I148").Value = "sd"
$ws.Range("J148").Value = "Statement-non-opinion"
$ws.Range("I151").Value = "sv"
$ws.Range("J151").Value = "Statement-opinion"
$ws.Range("I171").Value = "%"
$ws.Range("J171").Value = "Uninterpretable"
$ws.Range("I173").Value = "sd"
$ws.Range("J173").Value = "Statement-non-opinion"
$ws.Range("I193").Value = "ba"
$ws.Range("J193").Value = "Appreciation"
$ws.Range("I194").Value = "b"
$ws.Range("J194").Value = "Acknowledge (Backchannel)"
$ws.Range("I200").Value = "ba"
$ws.Range("J200").Value = "Appreciation"
$ws.Range("I216").Value = "sv"
$ws.Range("J216").Value = "Statement-opinion"
$ws.Range("I222").Value = "b"
$ws.Range("J222").Value = "Acknowledge (Backchannel)"
$ws.Range("I226").Value = "sd"
$ws.Range("J226").Value = "Statement-non-opinion"
$ws.Range("I230").Value = "sd"
$ws.Range("J230").Value = "Statement-non-opinion"
$ws.Range("I243").Value = "b"
$ws.Range("J243").Value = "Acknowledge (Backchannel)"
$ws.Range("I250").Value = "sv"
$ws.Range("J250").Value = "Statement-opinion"
$ws.Range("I252").Value = "%"
$ws.Range("J252").Value = "Uninterpretable"
$ws.Range("I256").Value = "ba"
$ws.Range("J256").Value = "Appreciation"
$ws.Range("I265").Value = "b"
$ws.Range("J265").Value = "Acknowledge (Backchannel)"
$ws.Range("I272").Value = "ba"
$ws.Range("J272").Value = "Appreciation"
$ws.Range("I285").Value = "b"
$ws.Range("J285").Value = "Acknowledge (Backchannel)"
$ws.Range("I292").Value = "sd"
$ws.Range("J292").Value = "Statement-non-opinion"
$ws.Range("I304").Value = "sv"
$ws.Range("J304").Value = "Statement-opinion"
$ws.Range("I309").Value = "ba"
$ws.Range("J309").Value = "Appreciation"
$ws.Range("I310").Value = "b"
$ws.Range("J310").Value = "Acknowledge (Backchannel)"
$ws.Range("I332").Value = "sv"
$ws.Range("J332").Value = "Statement-opinion"
$ws.Range("I349").Value = "sd"
$ws.Range("J349").Value = "Statement-non-opinion"
$ws.Range("I357").Value = "sv"
$ws.Range("J357").Value = "Statement-opinion"
$ws.Range("I365").Value = "sd"
$ws.Range("J365").Value = "Statement-non-opinion"
$ws.Range("I367").Value = "sd"
$ws.Range("J367").Value = "Statement-non-opinion"
$ws.Range("I376").Value = "b"
$ws.Range("J376").Value = "Acknowledge (Backchannel)"
$ws.Range("I386").Value = "aa"
$ws.Range("J386").Value = "Agree/Accept"
$ws.Range("I390").Value = "%"
$ws.Range("J390").Value = "Uninterpretable"
$ws.Range("I392").Value = "b"
$ws.Range("J392").Value = "Acknowledge (Backchannel)"
$ws.Range("I415").Value = "b"
$ws.Range("J415").Value = "Acknowledge (Backchannel)"
$ws.Range("I434").Value = "b"
$ws.Range("J434").Value = "Acknowledge (Backchannel)"
$ws.Range("I435").Value = "sd"
$ws.Range("J435").Value = "Statement-non-opinion"
$ws.Range("I446").Value = "b"
$ws.Range("J446").Value = "Acknowledge (Backchannel)"
$ws.Range("I455").Value = "sv"
$ws.Range("J455").Value = "Statement-opinion"
$ws.Range("I475").Value = "sd"
$ws.Range("J475").Value = "Statement-non-opinion"
$ws.Range("I477").Value = "b"
$ws.Range("J477").Value = "Acknowledge (Backchannel)"
$ws.Range("I480").Value = "sd"
$ws.Range("J480").Value = "Statement-non-opinion"
$ws.Range("I484").Value = "sv"
$ws.Range("J484").Value = "Statement-opinion"
$ws.Range("I485").Value = "%"
$ws.Range("J485").Value = "Uninterpretable"
$ws.Range("I486").Value = "aa"
$ws.Range("J486").Value = "Agree/Accept"
